$wb = $excel.ActiveWorkbook

# --- 1) Rename header cells (B1) on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header styling (bold, centered, bordered) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$data = @(
    @(2, 45004.99999999999, 140, -153.3884467673877, 412.0987444539376),
    @(3, 45011.99999999999, 149, -160.0266607614241, 444.7347502635249),
    @(4, 45025.99999999999, 166, -151.1660102993347, 447.1119777029601),
    @(5, 45032.99999999999, 175, -122.0913740386571, 474.3902664384001),
    @(6, 45039.99999999999, 183, -112.8724154280573, 484.3293208095423),
    @(7, 45046.99999999999, 192, -99.79959250736941, 489.7106316140446),
    @(8, 45053.99999999999, 201, -100.6900456015445, 510.0881398131596),
    @(9, 45060.99999999999, 209, -118.7491203365746, 498.090509550519),
    @(10, 45067.99999999999, 218, -82.77355822175892, 520.9055483810445),
    @(11, 45088.99999999999, 244, -60.48643904967059, 532.228930376116),
    @(12, 45172.99999999999, 348, 52.84546254613629, 648.9908924138175),
    @(13, 45179.99999999999, 357, 52.87865428938065, 665.2998176460627),
    @(14, 45214.99999999999, 400, 119.1866193518236, 710.2326131619543),
    @(15, 45221.99999999999, 409, 109.2170446373176, 720.713602856025),
    @(16, 45333.99999999999, 547, 248.5790042344327, 848.7581835239511),
    @(17, 45340.99999999999, 556, 263.3367641080357, 830.1641631059314),
    @(18, 45347.99999999999, 565, 234.8486650456693, 853.99577809579),
    @(19, 45354.99999999999, 573, 287.3063614169737, 877.0679714980828),
    @(20, 45361.99999999999, 582, 276.9490794012031, 886.5989614770468),
    @(21, 45368.99999999999, 591, 293.0188031708338, 877.5123435561206),
    @(22, 45375.99999999999, 599, 301.0780986087168, 882.6644955504717),
    @(23, 45382.99999999999, 608, 290.8928385568162, 926.4469041184445),
    @(24, 45389.99999999999, 617, 322.6137255101057, 901.2332794530615),
    @(25, 45396.99999999999, 625, 326.8126194227149, 912.8128985633957),
    @(26, 45403.99999999999, 634, 324.3606913952706, 941.6352667593054),
    @(27, 45410.99999999999, 643, 362.6908285357517, 954.1463627786086),
)

foreach ($row in $data) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}

# Match the date-style (format) used for the date column on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A27").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the originally active sheet/selection as it was before the edit
$wsWeekly.Activate()
